$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F5").Value = 7939
$ws1.Range("F9").Value = 6806
$ws1.Range("F17").Value = 88
$ws1.Range("F20").Value = 65
$ws1.Range("F25").Value = 3914
$ws1.Range("F27").Value = 382
$ws1.Range("F30").Value = 1533
$ws1.Range("F33").Value = 2814
$ws1.Range("F34").Value = 1977
$ws1.Range("F39").Value = 3828
$ws1.Range("F45").Value = 19
$ws1.Range("F46").Value = 1474

$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F17").Value = 244

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F7").Value = 7939
$ws4.Range("F9").Value = 6806
$ws4.Range("F15").Value = 88
$ws4.Range("F18").Value = 65
$ws4.Range("F22").Value = 3914
$ws4.Range("F26").Value = 382
$ws4.Range("F29").Value = 1533
$ws4.Range("F32").Value = 2814
$ws4.Range("F33").Value = 1977
$ws4.Range("F39").Value = 3828
$ws4.Range("F45").Value = 244
$ws4.Range("F46").Value = 1474
